$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the crypto rows that changed

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.123.19'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.03%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.652.90'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.04%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.21%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '218.64'
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.5220'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.20%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2653'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.22%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06350'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.89%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '20.43'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.55%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07697'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.27%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '4.631'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +3.41%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.674.29'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.47%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '1.882.19'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.15%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.5600'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.96%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.0₅8181'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +2.13%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '65.49'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.13%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '26.121.30'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("E19").Value = '  -0.14%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.638'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.11%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '10.47'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +3.86%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '191.97'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.38%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.942'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.21%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.24%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '145.50'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.03%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.1195'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.80%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.239'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.81%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '15.93'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.21%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.512'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.85%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.05470'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -4.09%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.271'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.24%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.455'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.82%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.370'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.30%  '

$ws.Range("E34").Value = '  -2.09%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.9529'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.21%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.787'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.28%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.400'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.44%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.5645'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.36%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01581'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.51%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '5.867'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.48%  '

$ws.Range("E41").Value = '  -0.17%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.8354'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.92%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.027.17'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.05%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '101.09'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.36%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.792.52'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.15%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '57.70'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.61%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.0₈108'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.81%  '

$ws.Range("E48").Value = '  -0.73%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.4341'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.25%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '7.972'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.02%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.05194'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -3.60%  '
